$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the MariaDB entries (B10 and B11) to match the new naming convention
$ws.Range("B10").Value = "MariaDB 10_6"
$ws.Range("B11").Value = "MariaDB 10_11"

# Update selection to B11 as the active cell (matches resulting selection in file)
$ws.Range("B11").Select()
